$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells being updated to be treated as text so that
# values like "1.00" or "0.998" keep their exact formatting instead of being
# auto-converted into numbers by Excel.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D13","D14","D15","D16","D18","D19","D20","D21","D22","D24","D28","D29","D30","D31","D32","D34","D36","D37","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "58.929.40"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.499.96"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "536.46"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "136.78"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.524.29"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "2.961.48"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "23.09"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "58.850.50"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.530.31"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "11.13"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "323.48"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "65.28"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "0.0₃0772"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "6.68"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "172.00"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "1.76"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "18.39"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "4.10"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "36.73"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "0.812"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "285.17"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "131.57"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.611"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.994"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "10.89"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "0.0922"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").Value = "0.0507"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "17.43"
$ws.Range("E51").Value = "  -2.21%  "

# Remove the temporary text formatting again so the cells keep their original
# (default/no explicit) style, matching the rest of the sheet.
foreach ($addr in $dCells) { $ws.Range($addr).ClearFormats() }

Write-Host "Applied 91 cell updates"
